# "updated main GSC export data"
# Append the new day's row (2025-11-18) to the end of the GSC HTTPS export table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$newRow = 44

# Write the date as a literal text value (matching how every other date in
# column A is stored - as text, not as a real Excel date serial). Using a
# formula that evaluates to the text, then pasting back as a value, avoids
# Excel's automatic "this looks like a date" reinterpretation that a plain
# .Value assignment of a date-shaped string would trigger.
$dateCell = $ws.Range("A" + $newRow)
$dateCell.Formula = '="2025-11-18"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("B" + $newRow).Value = 0
$ws.Range("C" + $newRow).Value = 29
